$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "08:00:16"
$ws.Range("D2").Value = "11:58:43"
$ws.Range("E2").Value = "12:58:08"
$ws.Range("F2").Value = "17:27:49"
$ws.Range("G2").Value = "08:28:08"
$ws.Range("C3").Value = "07:54:02"
$ws.Range("D3").Value = "12:18:26"
$ws.Range("E3").Value = "13:21:26"
$ws.Range("F3").Value = "17:31:39"
$ws.Range("G3").Value = "08:34:37"
$ws.Range("C4").Value = "07:38:42"
$ws.Range("D4").Value = "10:59:53"
$ws.Range("E4").Value = "11:56:57"
$ws.Range("F4").Value = "17:13:22"
$ws.Range("G4").Value = "08:37:36"
$ws.Range("C5").Value = "07:45:11"
$ws.Range("D5").Value = "10:48:22"
$ws.Range("E5").Value = "11:45:46"
$ws.Range("F5").Value = "17:14:53"
$ws.Range("G5").Value = "08:32:18"
$ws.Range("C6").Value = "08:20:04"
$ws.Range("D6").Value = "11:55:11"
$ws.Range("E6").Value = "12:54:49"
$ws.Range("F6").Value = "17:09:00"
$ws.Range("G6").Value = "07:49:18"
$ws.Range("C7").Value = "07:44:18"
$ws.Range("D7").Value = "11:57:22"
$ws.Range("E7").Value = "12:59:15"
$ws.Range("F7").Value = "17:10:39"
$ws.Range("G7").Value = "08:24:28"
$ws.Range("C8").Value = "07:53:53"
$ws.Range("D8").Value = "11:28:55"
$ws.Range("E8").Value = "12:28:44"
$ws.Range("F8").Value = "17:14:01"
$ws.Range("G8").Value = "08:20:19"
$ws.Range("C9").Value = "08:03:17"
$ws.Range("D9").Value = "11:35:20"
$ws.Range("E9").Value = "12:38:36"
$ws.Range("F9").Value = "17:05:03"
$ws.Range("G9").Value = "07:58:30"
$ws.Range("C10").Value = "08:09:52"
$ws.Range("D10").Value = "11:13:59"
$ws.Range("E10").Value = "12:15:31"
$ws.Range("F10").Value = "17:06:03"
$ws.Range("G10").Value = "07:54:39"
$ws.Range("C11").Value = "07:45:19"
$ws.Range("D11").Value = "11:51:59"
$ws.Range("E11").Value = "12:55:19"
$ws.Range("F11").Value = "17:05:00"
$ws.Range("G11").Value = "08:16:21"
$ws.Range("C12").Value = "07:47:59"
$ws.Range("D12").Value = "12:06:28"
$ws.Range("E12").Value = "12:56:45"
$ws.Range("F12").Value = "17:17:02"
$ws.Range("G12").Value = "08:38:46"
$ws.Range("C13").Value = "07:42:01"
$ws.Range("D13").Value = "11:14:34"
$ws.Range("E13").Value = "12:13:59"
$ws.Range("F13").Value = "17:30:04"
$ws.Range("G13").Value = "08:48:38"
$ws.Range("C14").Value = "07:59:57"
$ws.Range("D14").Value = "11:48:40"
$ws.Range("E14").Value = "12:48:40"
$ws.Range("F14").Value = "17:32:47"
$ws.Range("G14").Value = "08:32:50"
$ws.Range("C15").Value = "07:33:39"
$ws.Range("D15").Value = "11:59:56"
$ws.Range("E15").Value = "12:57:46"
$ws.Range("F15").Value = "17:33:40"
$ws.Range("G15").Value = "09:02:11"
$ws.Range("H15").Value = "00:12:11"
$ws.Range("C16").Value = "07:38:33"
$ws.Range("D16").Value = "11:27:13"
$ws.Range("E16").Value = "12:26:27"
$ws.Range("F16").Value = "17:31:35"
$ws.Range("G16").Value = "08:53:48"
$ws.Range("H16").Value = "00:03:48"
$ws.Range("C17").Value = "07:42:53"
$ws.Range("D17").Value = "11:39:24"
$ws.Range("E17").Value = "12:37:24"
$ws.Range("F17").Value = "17:31:52"
$ws.Range("G17").Value = "08:50:59"
$ws.Range("H17").Value = "00:00:59"
$ws.Range("C18").Value = "08:06:11"
$ws.Range("D18").Value = "12:06:11"
$ws.Range("E18").Value = "13:06:11"
$ws.Range("F18").Value = "17:42:42"
$ws.Range("G18").Value = "08:36:31"
$ws.Range("G19").Value = "144:19:57"
$ws.Range("H19").Value = "00:16:58"
